$d = $word.ActiveDocument

# The original paragraph ends with "Test 1" followed immediately (with no
# visible text in between) by the hidden "_GoBack" bookmark. We want to
# split the paragraph right after "Test 1" so that a brand-new, bold
# paragraph containing "Test 2" is created, with the "_GoBack" bookmark
# now wrapping the new "Test 2" run instead of sitting at the end of the
# first paragraph.

# Locate "Test 1" and collapse the range to just after it.
$findRange = $d.Content
$found = $findRange.Find.Execute("Test 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Collapse(0)  # wdCollapseEnd
$splitPos = $findRange.End

# Remove the pre-existing hidden bookmark so we can recreate it in the
# right place afterwards (it is hidden from Bookmarks.Count, but is still
# reachable by name).
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to clean up
}

# Split "Test 1" into its own paragraph, leaving a new, empty paragraph
# right after it.
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# Figure out where the new (now second) paragraph begins and type
# "Test 2" into it.
$newText = "Test 2"
$newParaStart = $splitPos + 1
$newParaRange = $d.Range($newParaStart, $newParaStart)
$newParaRange.InsertAfter($newText)
$newParaEnd = $newParaStart + $newText.Length

# Re-create the "_GoBack" bookmark so it wraps the freshly-typed "Test 2".
$d.Bookmarks.Add("_GoBack", $d.Range($newParaStart, $newParaEnd))

# Bold the whole new paragraph (run text + paragraph mark), matching a
# paragraph-level bold toggle applied while the cursor/selection was in
# that paragraph.
$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Bold = 1
